# ---------------------------------------------------------------------------
# CasosColombia.xlsx update
#   1. Six previously-numeric daily counts are recorded as "NaN" (no data
#      reported that day) instead of a number.
#   2. Two new daily rows (2020-08-31 and 2020-09-01) are appended at the
#      bottom of the table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Cells that flip from a numeric reading to the literal text "NaN" ---
$nanCells = "DT9","L19","AF127","BQ137","BQ138","W167"
foreach ($addr in $nanCells) {
    $ws.Range($addr).Value = "NaN"
}

# --- 2. Append the two new rows (columns A..DX = 128 columns each) ---
$row180 = @(44074,615168,2701,80576,64127,211300,25267,3612,2964,6067,5285,10815,3722,19738,22925,5139,4920,12461,7990,14086,11850,2895,1276,5992,18162,11709,7129,46944,1062,219,340,447,140,95,272,1957,3050,36130,6717,2402,36864,916,20443,1460,8278,1517,1558,4353,1611,935,2469,2595,47050,12226,2745,7524,3926,279,1398,2585,729,2020,8415,8361,8471,13804,1882,823,7750,6684,7868,1546,1523,3027,3235,902,4275,2434,1312,700,2096,1878,1235,954,4965,1384,1163,1260,1592,1483,1676,1135,1070,1097,610,3005,1020,804,742,1287,1172,647,728,918,1176,991,1123,881,317,333,678,584,404,530,328,596,702,512,475,363,512,117809,259597,10505,112013,70037,30061,9102)
$row181 = @(44075,624069,2704,81489,64255,213631,25454,3701,3022,6165,5406,11211,3742,19953,23609,5266,5113,12529,8354,14337,12047,3005,1304,6267,19039,11778,7358,47658,1088,262,340,448,156,114,279,1964,3179,36218,6785,2403,37428,918,20583,1462,8409,1531,1559,4539,1629,936,2472,2596,47683,12321,2861,7704,4143,279,1398,2602,730,2028,8464,8401,8622,13823,1889,825,8018,6987,8336,1572,1547,3101,3363,928,4324,2463,1346,727,2146,1896,1276,967,5041,1468,1182,1304,1644,1539,1746,1143,1083,1099,620,3007,1056,809,751,1314,1230,661,735,942,1180,999,1148,886,317,333,681,588,408,530,330,605,703,513,475,363,512,118715,263487,10866,113210,71143,31135,9326)

$newRowNum = 180
foreach ($rowVals in @($row180, $row181)) {
    $colIdx = 1
    foreach ($v in $rowVals) {
        $ws.Cells.Item($newRowNum, $colIdx).Value = $v
        $colIdx = $colIdx + 1
    }
    $newRowNum = $newRowNum + 1
}

# --- 3. Leave the selection where the edit finished (last row, ColREje col) ---
$ws.Range("DW180").Select()
